$wb = $excel.ActiveWorkbook

# Sheet "Auth": update the B2:B6 timestamp values (column "Date")
$wsAuth = $wb.Worksheets.Item("Auth")
$wsAuth.Range("B2").Value = "Fri Oct 21 13:51:30 EDT 2022"
$wsAuth.Range("B3").Value = "Fri Oct 21 13:51:46 EDT 2022"
$wsAuth.Range("B4").Value = "Fri Oct 21 13:52:01 EDT 2022"
$wsAuth.Range("B5").Value = "Fri Oct 21 13:52:15 EDT 2022"
$wsAuth.Range("B6").Value = "Fri Oct 21 13:52:30 EDT 2022"

# Sheet "AuthCapture": update the B2:B6 timestamp values (column "Date")
$wsAuthCapture = $wb.Worksheets.Item("AuthCapture")
$wsAuthCapture.Range("B2").Value = "Fri Oct 21 13:52:47 EDT 2022"
$wsAuthCapture.Range("B3").Value = "Fri Oct 21 13:53:11 EDT 2022"
$wsAuthCapture.Range("B4").Value = "Fri Oct 21 13:53:34 EDT 2022"
$wsAuthCapture.Range("B5").Value = "Fri Oct 21 13:53:57 EDT 2022"
$wsAuthCapture.Range("B6").Value = "Fri Oct 21 13:54:23 EDT 2022"
